$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce precision (custom accuracy) of row 5 values from 3 decimals to 2 decimals
$ws.Range("B5").Value = 11.11
$ws.Range("C5").Value = 8.58
$ws.Range("D5").Value = 0.27
$ws.Range("E5").Value = 24.74
$ws.Range("F5").Value = 20.36
$ws.Range("G5").Value = 8.95
$ws.Range("H5").Value = 35.27
$ws.Range("I5").Value = 14.02
$ws.Range("J5").Value = 6.42
$ws.Range("K5").Value = 9.65
$ws.Range("L5").Value = 10.28
$ws.Range("M5").Value = 10.56
$ws.Range("N5").Value = 2.9
$ws.Range("O5").Value = 8.8
$ws.Range("P5").Value = 12.83
$ws.Range("Q5").Value = 7.49
$ws.Range("R5").Value = 0.16
$ws.Range("S5").Value = 0.42
$ws.Range("T5").Value = 129.86
$ws.Range("U5").Value = 25.25
$ws.Range("V5").Value = 8.22
$ws.Range("W5").Value = 16.79
$ws.Range("X5").Value = 9.29
$ws.Range("Y5").Value = 1.18
$ws.Range("Z5").Value = 17.26
$ws.Range("AA5").Value = 7.34
$ws.Range("AB5").Value = 6.87
$ws.Range("AC5").Value = 8.07
$ws.Range("AD5").Value = 10.81
$ws.Range("AE5").Value = 0.21
$ws.Range("AF5").Value = 32.07
$ws.Range("AG5").Value = 4.94
$ws.Range("AH5").Value = 10.24

# Remove the last data row (row 6), shrinking the dataset
$ws.Rows.Item(6).Delete()

# Narrow a couple of columns (W and AD) to width 7
$ws.Range("W1").ColumnWidth = 6.1699999999999999
$ws.Range("AD1").ColumnWidth = 6.1699999999999999
